$d = $word.ActiveDocument

$pairs = @(
    @("252×3=756", "552×7=3864"),
    @("820×6=4920", "259×3=777"),
    @("888×6=5328", "140×8=1120"),
    @("230×4=920", "438×2=876"),
    @("933×5=4665", "562×9=5058"),
    @("979×6=5874", "678×4=2712"),
    @("540×9=4860", "312×6=1872"),
    @("416×8=3328", "842×4=3368"),
    @("199×3=597", "104×9=936"),
    @("734×5=3670", "328×5=1640"),
    @("994×2=1988", "528×9=4752"),
    @("566×7=3962", "627×6=3762"),
    @("499×2=998", "520×4=2080"),
    @("610×4=2440", "874×6=5244"),
    @("463×8=3704", "189×6=1134"),
    @("359×9=3231", "890×5=4450"),
    @("979×5=4895", "332×3=996"),
    @("896×8=7168", "411×6=2466"),
    @("927×5=4635", "775×2=1550"),
    @("638×7=4466", "631×6=3786"),
    @("561×7=3927", "415×2=830"),
    @("578×4=2312", "616×3=1848"),
    @("444×5=2220", "783×7=5481"),
    @("674×5=3370", "416×6=2496"),
    @("937×3=2811", "850×5=4250")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
